$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, $val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.NumberFormat = "General"
    $rng.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "246.18"
Set-TextValue $ws.Range("E2") "0.26%"
Set-TextValue $ws.Range("D3") "29.74"
Set-TextValue $ws.Range("E3") "-0.55%"
Set-TextValue $ws.Range("D4") "5.165"
Set-TextValue $ws.Range("E4") "0.66%"
Set-TextValue $ws.Range("E5") "1.17%"
Set-TextValue $ws.Range("D6") "6.682"
Set-TextValue $ws.Range("E6") "1.80%"
Set-TextValue $ws.Range("D7") "3.225"
Set-TextValue $ws.Range("E7") "6.97%"
Set-TextValue $ws.Range("E8") "-0.42%"
Set-TextValue $ws.Range("D9") "0.8623"
Set-TextValue $ws.Range("E9") "-1.02%"
Set-TextValue $ws.Range("D10") "0.1376"
Set-TextValue $ws.Range("E10") "2.21%"
Set-TextValue $ws.Range("D11") "0.07103"
Set-TextValue $ws.Range("E11") "2.78%"
Set-TextValue $ws.Range("D12") "0.03204"
Set-TextValue $ws.Range("E12") "11.04%"
Set-TextValue $ws.Range("D13") "0.09377"
Set-TextValue $ws.Range("E13") "0.10%"
Set-TextValue $ws.Range("D14") "0.001537"
Set-TextValue $ws.Range("E14") "0.86%"
Set-TextValue $ws.Range("D15") "0.0005980"
Set-TextValue $ws.Range("E15") "-94.09%"
Set-TextValue $ws.Range("D16") "0.005891"
Set-TextValue $ws.Range("E16") "-1.40%"
Set-TextValue $ws.Range("D17") "3.499"
Set-TextValue $ws.Range("E17") "-0.29%"
Set-TextValue $ws.Range("E18") "1.49%"
Set-TextValue $ws.Range("D20") "0.03346"
Set-TextValue $ws.Range("E20") "0.87%"
Set-TextValue $ws.Range("D21") "0.1300"
Set-TextValue $ws.Range("E21") "-0.34%"
Set-TextValue $ws.Range("D22") "3.491"
Set-TextValue $ws.Range("E22") "-2.71%"
Set-TextValue $ws.Range("D23") "0.04144"
Set-TextValue $ws.Range("E24") "0.52%"
Set-TextValue $ws.Range("D25") "0.001227"
Set-TextValue $ws.Range("E25") "1.41%"
Set-TextValue $ws.Range("D26") "0.004143"
Set-TextValue $ws.Range("E26") "-7.58%"
Set-TextValue $ws.Range("E27") "2.06%"
Set-TextValue $ws.Range("E28") "4.27%"
Set-TextValue $ws.Range("D40") "0.03761"
Set-TextValue $ws.Range("E40") "-0.31%"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue $ws.Range("D41") "0.005678"
Set-TextValue $ws.Range("E41") "-0.52%"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws.Range("D42") "0.1070"
Set-TextValue $ws.Range("E42") "0.33%"
Set-TextValue $ws.Range("D43") "0.002199"
Set-TextValue $ws.Range("E43") "-0.82%"
Set-TextValue $ws.Range("D44") "0.009552"
Set-TextValue $ws.Range("E44") "-0.44%"
Set-TextValue $ws.Range("D45") "0.00005296"
Set-TextValue $ws.Range("E45") "4.57%"
Set-TextValue $ws.Range("E46") "0.34%"
Set-TextValue $ws.Range("D47") "0.05801"
Set-TextValue $ws.Range("E47") "-27.22%"
Set-TextValue $ws.Range("D48") "0.002182"
Set-TextValue $ws.Range("E48") "-20.49%"
Set-TextValue $ws.Range("D49") "0.00002099"
Set-TextValue $ws.Range("E49") "0.34%"
Set-TextValue $ws.Range("D50") "0.0001999"
Set-TextValue $ws.Range("E50") "0.34%"
